$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain decimal number need to be forced to Text
# format first, otherwise Excel auto-converts the numeric-looking string to a Number
# (dropping the formatted trailing zero), same as typing it into a General cell would.

$ws.Range("D2").Value = "43.982.12"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "2.242.52"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "270.62"
$ws.Range("E5").Value = "  +4.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.80"
$ws.Range("E6").Value = "  +15.84%  "
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +6.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.97"
$ws.Range("E10").Value = "  +6.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0966"
$ws.Range("E11").Value = "  +5.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.28"
$ws.Range("E12").Value = "  +18.99%  "
$ws.Range("D14").Value = "2.580.68"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.17"
$ws.Range("E15").Value = "  +6.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.814"
$ws.Range("E16").Value = "  +4.75%  "
$ws.Range("D17").Value = "2.249.56"
$ws.Range("E17").Value = "  +3.28%  "
$ws.Range("D18").Value = "43.933.63"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("E20").Value = "  +3.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.74"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.33"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.93"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.44"
$ws.Range("E26").Value = "  +7.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.51"
$ws.Range("E27").Value = "  +12.19%  "
$ws.Range("E28").Value = "  +6.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.73"
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.82"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  +5.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.98"
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.49"
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0354"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.34"
$ws.Range("E38").Value = "  -2.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.62"
$ws.Range("E39").Value = "  +27.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.93"
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.227"
$ws.Range("E41").Value = "  +14.23%  "
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.65"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0998"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.94"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.39"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("E48").Value = "  +4.49%  "
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("D51").Value = "2.461.07"
$ws.Range("E51").Value = "  +1.94%  "
